# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds the "K" values. Update the computed K values for rows 2-9.
$kValues = @{
    2 = 0
    3 = 0
    4 = 1
    5 = 1
    6 = 1
    7 = 0
    8 = 0
    9 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
